$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume (E) columns to retain text formatting so numeric-looking
# strings (e.g. "53.07") are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.429.24"
$ws.Range("E2").Value = "  -3.66%  "
$ws.Range("D3").Value = "1.993.74"
$ws.Range("E3").Value = "  -6.22%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "329.28"
$ws.Range("E5").Value = "  -5.31%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "0.5002"
$ws.Range("E7").Value = "  -4.53%  "
$ws.Range("D8").Value = "0.4222"
$ws.Range("E8").Value = "  -5.85%  "
$ws.Range("D9").Value = "53.07"
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("D10").Value = "0.08907"
$ws.Range("E10").Value = "  -5.33%  "
$ws.Range("D11").Value = "1.120"
$ws.Range("E11").Value = "  -5.52%  "
$ws.Range("D12").Value = "23.16"
$ws.Range("E12").Value = "  -8.51%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "8.069"
$ws.Range("E13").Value = "  -7.61%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.997.64"
$ws.Range("E14").Value = "  -5.07%  "
$ws.Range("E15").Value = "  -6.91%  "
$ws.Range("D16").Value = "95.81"
$ws.Range("E16").Value = "  -6.54%  "
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "0.00001108"
$ws.Range("E18").Value = "  -5.20%  "
$ws.Range("D19").Value = "0.06614"
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("D20").Value = "19.65"
$ws.Range("E20").Value = "  -8.77%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "5.956"
$ws.Range("E22").Value = "  -6.25%  "
$ws.Range("D23").Value = "29.463.02"
$ws.Range("E23").Value = "  -3.60%  "
$ws.Range("E24").Value = "  -7.26%  "
$ws.Range("D25").Value = "2.259"
$ws.Range("E25").Value = "  -3.25%  "
$ws.Range("D26").Value = "158.19"
$ws.Range("E26").Value = "  -2.79%  "
$ws.Range("D27").Value = "20.61"
$ws.Range("E27").Value = "  -7.51%  "
$ws.Range("D28").Value = "6.476"
$ws.Range("E28").Value = "  -6.71%  "
$ws.Range("D29").Value = "2.323"
$ws.Range("E29").Value = "  -9.15%  "
$ws.Range("D30").Value = "127.67"
$ws.Range("E30").Value = "  -5.19%  "
$ws.Range("E31").Value = "  -10.10%  "
$ws.Range("D32").Value = "0.09932"
$ws.Range("E32").Value = "  -6.57%  "
$ws.Range("D33").Value = "1.565"
$ws.Range("E33").Value = "  -12.33%  "
$ws.Range("D34").Value = "5.838"
$ws.Range("E34").Value = "  -7.49%  "
$ws.Range("D35").Value = "3.783"
$ws.Range("E35").Value = "  -4.65%  "
$ws.Range("D36").Value = "9.578"
$ws.Range("E36").Value = "  -10.28%  "
$ws.Range("D37").Value = "0.02459"
$ws.Range("E37").Value = "  -7.73%  "
$ws.Range("D38").Value = "0.06331"
$ws.Range("E38").Value = "  -8.03%  "
$ws.Range("E39").Value = "  -3.81%  "
$ws.Range("D40").Value = "0.6512"
$ws.Range("E40").Value = "  -8.89%  "
$ws.Range("D41").Value = "11.67"
$ws.Range("E41").Value = "  -8.10%  "
$ws.Range("D42").Value = "0.2064"
$ws.Range("E42").Value = "  -8.49%  "
$ws.Range("D44").Value = "0.6325"
$ws.Range("E44").Value = "  -8.94%  "
$ws.Range("D45").Value = "13.46"
$ws.Range("E45").Value = "  -7.74%  "
$ws.Range("E46").Value = "  -8.44%  "
$ws.Range("D47").Value = "1.286"
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("E48").Value = "  -3.73%  "
$ws.Range("D49").Value = "0.00000000326"
$ws.Range("E49").Value = "  -5.14%  "
$ws.Range("D50").Value = "0.06990"
$ws.Range("E50").Value = "  -3.42%  "
$ws.Range("D51").Value = "1.134"
$ws.Range("E51").Value = "  -5.43%  "

# Restore the original General number format now that the text values are locked in.
$ws.Range("D2:E51").NumberFormat = "General"
